$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new data rows (67-77) to Sheet1, mirroring the upstream CSV upload.
# Columns: A=z, D=xHI, E=dxHI_p_68, F=dxHI_m_68, G=dxHI_p_95, H=dxHI_m_95, I=Method, J=Reference

# --- Numeric / value columns (A, D, E, F, G, H) ---
$ws.Range("A67").Value = 7.3
$ws.Range("D67").Value = 0.75
$ws.Range("E67").Value = 0.09
$ws.Range("F67").Value = 0.13
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 0

$ws.Range("A68").Value = 7
$ws.Range("D68").Value = 0.18
$ws.Range("E68").Value = 0.14
$ws.Range("F68").Value = 0.12
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0

$ws.Range("A69").Value = 6.6
$ws.Range("D69").Value = 0.21
$ws.Range("E69").Value = 0.19
$ws.Range("F69").Value = 0.14
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0

$ws.Range("A70").Value = 6.6
$ws.Range("D70").Value = 0.15
$ws.Range("E70").Value = 0.1
$ws.Range("F70").Value = 0.08
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0

$ws.Range("A71").Value = 5.7
$ws.Range("D71").Value = 0.06
$ws.Range("E71").Value = 0.12
$ws.Range("F71").Value = 0.03
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0

$ws.Range("A72").Value = 5.7
$ws.Range("D72").Value = "<0.05"
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 0

$ws.Range("A73").Value = 10.4
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0.4
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0

$ws.Range("A74").Value = 8.6
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0.2
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0

$ws.Range("A75").Value = 7
$ws.Range("D75").Value = 0.65
$ws.Range("E75").Value = 0.27
$ws.Range("F75").Value = 0.35
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0

$ws.Range("A76").Value = 5.8
$ws.Range("D76").Value = 0.25
$ws.Range("E76").Value = 0.1
$ws.Range("F76").Value = 0.2
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 0

$ws.Range("A77").Value = 5
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0.12
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0

# --- Method column (I), filled bottom-up ---
$ws.Range("I77").Value = "GP trough, Damping wing"
$ws.Range("I76").Value = "GP trough, Damping wing"
$ws.Range("I75").Value = "GP trough, Damping wing"
$ws.Range("I74").Value = "GP trough, Damping wing"
$ws.Range("I73").Value = "GP trough, Damping wing"
$ws.Range("I72").Value = "Lya LF"
$ws.Range("I71").Value = "Angular Correlation function"
$ws.Range("I70").Value = "Lya LF"
$ws.Range("I69").Value = "Angular Correlation function"
$ws.Range("I68").Value = "Lya LF"
$ws.Range("I67").Value = "Lya LF"

# --- Reference column (J), filled top-down ---
$ws.Range("J67").Value = "Umeda+25a"
$ws.Range("J68").Value = "Umeda+25a"
$ws.Range("J69").Value = "Umeda+25a"
$ws.Range("J70").Value = "Umeda+25a"
$ws.Range("J71").Value = "Umeda+25a"
$ws.Range("J72").Value = "Umeda+25a"
$ws.Range("J73").Value = "Umeda+25b"
$ws.Range("J74").Value = "Umeda+25b"
$ws.Range("J75").Value = "Umeda+25b"
$ws.Range("J76").Value = "Umeda+25b"
$ws.Range("J77").Value = "Umeda+25b"

# Reflect the author's final on-screen selection/viewport from the diff.
[void]$ws.Range("D61:F64").Select()
[void]$ws.Range("H85").Select()
